$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Addr, $Val) {
    $cell = $Sheet.Range($Addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws 'D2' '25.946.14'
Set-TextValue $ws 'E2' '  +0.45%  '
# Row 3
Set-TextValue $ws 'D3' '1.645.51'
Set-TextValue $ws 'E3' '  +0.74%  '
# Row 4
Set-TextValue $ws 'D4' '1.006'
Set-TextValue $ws 'E4' '  +0.43%  '
# Row 5
Set-TextValue $ws 'D5' '215.77'
Set-TextValue $ws 'E5' '  +0.45%  '
# Row 6
Set-TextValue $ws 'D6' '0.5108'
Set-TextValue $ws 'E6' '  +1.90%  '
# Row 7
Set-TextValue $ws 'D7' '1.005'
Set-TextValue $ws 'E7' '  +0.32%  '
# Row 8
Set-TextValue $ws 'D8' '0.2582'
Set-TextValue $ws 'E8' '  +0.76%  '
# Row 9
Set-TextValue $ws 'D9' '0.06423'
Set-TextValue $ws 'E9' '  +0.50%  '
# Row 10
Set-TextValue $ws 'D10' '19.68'
Set-TextValue $ws 'E10' '  +0.56%  '
# Row 11
Set-TextValue $ws 'D11' '0.07783'
Set-TextValue $ws 'E11' '  +1.40%  '
# Row 12
Set-TextValue $ws 'D12' '4.320'
Set-TextValue $ws 'E12' '  +1.97%  '
# Row 13
Set-TextValue $ws 'D13' '1.644.21'
Set-TextValue $ws 'E13' '  +0.34%  '
# Row 14
Set-TextValue $ws 'D14' '0.5470'
Set-TextValue $ws 'E14' '  +0.98%  '
# Row 15
Set-TextValue $ws 'D15' '0.0₅7895'
Set-TextValue $ws 'E15' '  -0.17%  '
# Row 16
Set-TextValue $ws 'D16' '64.84'
Set-TextValue $ws 'E16' '  +2.20%  '
# Row 17
Set-TextValue $ws 'D17' '26.012.58'
Set-TextValue $ws 'E17' '  +0.65%  '
# Row 18
Set-TextValue $ws 'D18' '1.005'
Set-TextValue $ws 'E18' '  +0.23%  '
# Row 19
Set-TextValue $ws 'D19' '198.60'
Set-TextValue $ws 'E19' '  -1.56%  '
# Row 20
Set-TextValue $ws 'D20' '4.465'
Set-TextValue $ws 'E20' '  +3.27%  '
# Row 21
Set-TextValue $ws 'D21' '10.03'
Set-TextValue $ws 'E21' '  +1.18%  '
# Row 22
Set-TextValue $ws 'D22' '6.078'
Set-TextValue $ws 'E22' '  +2.03%  '
# Row 23
Set-TextValue $ws 'E23' '  +0.50%  '
# Row 24
Set-TextValue $ws 'D24' '1.860'
Set-TextValue $ws 'E24' '  -2.71%  '
# Row 25
Set-TextValue $ws 'D25' '140.18'
Set-TextValue $ws 'E25' '  -0.76%  '
# Row 26
Set-TextValue $ws 'D26' '0.1151'
Set-TextValue $ws 'E26' '  +1.22%  '
# Row 27
Set-TextValue $ws 'D27' '6.905'
Set-TextValue $ws 'E27' '  +3.17%  '
# Row 28
Set-TextValue $ws 'D28' '15.77'
Set-TextValue $ws 'E28' '  +0.61%  '
# Row 29
Set-TextValue $ws 'D29' '1.241'
Set-TextValue $ws 'E29' '  +0.16%  '
# Row 30
Set-TextValue $ws 'D30' '0.05029'
Set-TextValue $ws 'E30' '  +1.01%  '
# Row 31
Set-TextValue $ws 'D31' '3.290'
Set-TextValue $ws 'E31' '  +1.02%  '
# Row 32
Set-TextValue $ws 'D32' '3.203'
Set-TextValue $ws 'E32' '  +1.00%  '
# Row 33
Set-TextValue $ws 'E33' '  +0.71%  '
# Row 34
Set-TextValue $ws 'D34' '2.361'
Set-TextValue $ws 'E34' '  -0.12%  '
# Row 35
Set-TextValue $ws 'D35' '0.8943'
Set-TextValue $ws 'E35' '  +0.50%  '
# Row 36
Set-TextValue $ws 'D36' '2.589'
Set-TextValue $ws 'E36' '  -0.90%  '
# Row 37
Set-TextValue $ws 'B37' 'Maker'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D37' '1.135.16'
Set-TextValue $ws 'E37' '  -2.93%  '
# Row 38
Set-TextValue $ws 'B38' 'ImmutableX'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D38' '0.5549'
Set-TextValue $ws 'E38' '  -0.45%  '
# Row 39
Set-TextValue $ws 'D39' '0.01562'
Set-TextValue $ws 'E39' '  +0.35%  '
# Row 40
Set-TextValue $ws 'E40' '  +0.40%  '
# Row 41
Set-TextValue $ws 'D41' '5.671'
Set-TextValue $ws 'E41' '  -0.11%  '
# Row 42
Set-TextValue $ws 'D42' '0.8173'
Set-TextValue $ws 'E42' '  +1.27%  '
# Row 43
Set-TextValue $ws 'E43' '  +9.99%  '
# Row 44
Set-TextValue $ws 'E44' '  +0.71%  '
# Row 45
Set-TextValue $ws 'D45' '1.786.34'
Set-TextValue $ws 'E45' '  +0.91%  '
# Row 46
Set-TextValue $ws 'D46' '0.4531'
Set-TextValue $ws 'E46' '  +0.36%  '
# Row 47
Set-TextValue $ws 'D47' '55.25'
Set-TextValue $ws 'E47' '  +1.17%  '
# Row 48
Set-TextValue $ws 'D48' '1.004'
Set-TextValue $ws 'E48' '  +0.20%  '
# Row 49
Set-TextValue $ws 'D49' '0.05092'
Set-TextValue $ws 'E49' '  +0.33%  '
# Row 50
Set-TextValue $ws 'D50' '0.09577'
Set-TextValue $ws 'E50' '  +3.63%  '
# Row 51
Set-TextValue $ws 'D51' '1.007'
Set-TextValue $ws 'E51' '  +0.31%  '
